$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 383
$ws.Range("F6").Value = 538
$ws.Range("F7").Value = 50
$ws.Range("F8").Value = 9673
$ws.Range("F10").Value = 2606
$ws.Range("G10").Value = 19.9
$ws.Range("F12").Value = 2380
$ws.Range("F13").Value = 2623
$ws.Range("F15").Value = 270
$ws.Range("F16").Value = 2050
$ws.Range("F18").Value = 74
$ws.Range("F19").Value = 361
$ws.Range("F22").Value = 295
$ws.Range("F23").Value = 58
$ws.Range("F24").Value = 131
$ws.Range("F26").Value = 1269
$ws.Range("F27").Value = 1237
$ws.Range("F28").Value = 90
$ws.Range("F29").Value = 118
$ws.Range("F30").Value = 244
$ws.Range("F31").Value = 1645
$ws.Range("F32").Value = 2745
$ws.Range("F34").Value = 975
$ws.Range("F35").Value = 342
$ws.Range("F37").Value = 38
$ws.Range("F38").Value = 46
$ws.Range("F40").Value = 23

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 7
$ws.Range("F14").Value = 148

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 941
$ws.Range("F4").Value = 110

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 941
$ws.Range("F5").Value = 110
$ws.Range("F7").Value = 383
$ws.Range("F8").Value = 7
$ws.Range("F10").Value = 538
$ws.Range("F11").Value = 50
$ws.Range("F12").Value = 9673
$ws.Range("F15").Value = 2606
$ws.Range("G15").Value = 19.9
$ws.Range("F17").Value = 2380
$ws.Range("F18").Value = 2623
$ws.Range("F20").Value = 270
$ws.Range("F21").Value = 2050
$ws.Range("F23").Value = 74
$ws.Range("F24").Value = 361
$ws.Range("F27").Value = 295
$ws.Range("F28").Value = 58
$ws.Range("F29").Value = 131
$ws.Range("F31").Value = 1269
$ws.Range("F32").Value = 1237
$ws.Range("F33").Value = 90
$ws.Range("F34").Value = 118
$ws.Range("F36").Value = 1645
$ws.Range("F38").Value = 2745
$ws.Range("F39").Value = 975
$ws.Range("F41").Value = 342
$ws.Range("F46").Value = 46
$ws.Range("F48").Value = 23
$ws.Range("F49").Value = 148
$ws.Range("F50").Value = 148
